$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.714.37"
$ws.Range("E2").Value = "  -6.99%  "
$ws.Range("D3").Value = "3.680.43"
$ws.Range("E3").Value = "  -6.86%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'582.18"
$ws.Range("E5").Value = "  -3.75%  "
$ws.Range("D6").Value = "'169.78"
$ws.Range("E6").Value = "  -1.86%  "
$ws.Range("D7").Value = "3.670.30"
$ws.Range("E7").Value = "  -6.96%  "
$ws.Range("D8").Value = "'0.621"
$ws.Range("E8").Value = "  -9.23%  "
$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "  +0.20%  "
$ws.Range("D10").Value = "'0.701"
$ws.Range("E10").Value = "  -11.42%  "
$ws.Range("E11").Value = "  -12.26%  "
$ws.Range("D12").Value = "'51.47"
$ws.Range("E12").Value = "  -8.51%  "
$ws.Range("D13").Value = "'0.0000286"
$ws.Range("E13").Value = "  -12.94%  "
$ws.Range("D14").Value = "'10.39"
$ws.Range("E14").Value = "  -10.68%  "
$ws.Range("D15").Value = "4.245.77"
$ws.Range("E15").Value = "  -7.29%  "
$ws.Range("D16").Value = "3.689.81"
$ws.Range("E16").Value = "  -6.64%  "
$ws.Range("D17").Value = "'19.27"
$ws.Range("E17").Value = "  -10.48%  "
$ws.Range("E18").Value = "  -3.41%  "
$ws.Range("D19").Value = "'12.74"
$ws.Range("E19").Value = "  -9.76%  "
$ws.Range("E20").Value = "  -9.79%  "
$ws.Range("D21").Value = "67.472.78"
$ws.Range("E21").Value = "  -7.20%  "
$ws.Range("D22").Value = "'403.02"
$ws.Range("E22").Value = "  -9.46%  "
$ws.Range("E23").Value = "  -6.93%  "
$ws.Range("D24").Value = "'87.53"
$ws.Range("E24").Value = "  -8.69%  "
$ws.Range("D25").Value = "'3.03"
$ws.Range("E25").Value = "  -9.24%  "
$ws.Range("D26").Value = "'12.69"
$ws.Range("E26").Value = "  -10.69%  "
$ws.Range("D27").Value = "'10.83"
$ws.Range("E27").Value = "  -3.85%  "
$ws.Range("D28").Value = "'5.98"
$ws.Range("E28").Value = "  +1.39%  "
$ws.Range("E29").Value = "  -11.09%  "
$ws.Range("D30").Value = "'9.40"
$ws.Range("E30").Value = "  -10.06%  "
$ws.Range("D31").Value = "'32.42"
$ws.Range("E31").Value = "  -9.87%  "
$ws.Range("D32").Value = "'7.47"
$ws.Range("E32").Value = "  -6.35%  "
$ws.Range("D33").Value = "'12.34"
$ws.Range("E33").Value = "  -11.57%  "
$ws.Range("E34").Value = "  -10.20%  "
$ws.Range("D35").Value = "'64.51"
$ws.Range("E35").Value = "  -6.98%  "
$ws.Range("D36").Value = "'43.15"
$ws.Range("E36").Value = "  -13.37%  "
$ws.Range("D37").Value = "'591.90"
$ws.Range("E37").Value = "  -6.58%  "
$ws.Range("D38").Value = "0.0₃0886"
$ws.Range("E38").Value = "  -11.20%  "
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  +0.11%  "
$ws.Range("D40").Value = "'0.394"
$ws.Range("E40").Value = "  -8.20%  "
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("D42").Value = "'0.134"
$ws.Range("E42").Value = "  -8.20%  "
$ws.Range("D43").Value = "'2.75"
$ws.Range("E43").Value = "  +4.13%  "
$ws.Range("D44").Value = "'2.97"
$ws.Range("E44").Value = "  -13.47%  "
$ws.Range("D45").Value = "'0.0432"
$ws.Range("E45").Value = "  -9.83%  "
$ws.Range("D46").Value = "'2.78"
$ws.Range("E46").Value = "  -14.25%  "
$ws.Range("D47").Value = "'9.15"
$ws.Range("E47").Value = "  -13.96%  "
$ws.Range("D48").Value = "2.745.64"
$ws.Range("E48").Value = "  -3.43%  "
$ws.Range("E49").Value = "  -10.59%  "
$ws.Range("D50").Value = "'3.14"
$ws.Range("E50").Value = "  -7.98%  "
$ws.Range("E51").Value = "  -3.93%  "
